$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - odds updated for Wilstermann vs Tomayapo
$ws.Range("G2").Value = 2.15
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 1.06
$ws.Range("K2").Value = 9.5
$ws.Range("L2").Value = 1.29
$ws.Range("M2").Value = 3.5
$ws.Range("N2").Value = 1.98
$ws.Range("O2").Value = 1.83
$ws.Range("P2").Value = 1.4
$ws.Range("Q2").Value = 2.75
$ws.Range("W2").Value = 19
$ws.Range("X2").Value = 17
$ws.Range("Y2").Value = 26
$ws.Range("Z2").Value = 9.5
$ws.Range("AA2").Value = 6
$ws.Range("AE2").Value = 11
$ws.Range("AF2").Value = 17

# Row 5 - odds newly populated for Kuressaare vs Flora
$ws.Range("G5").Value = 7.2
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 1.3
$ws.Range("T5").Value = 20
$ws.Range("U5").Value = 40
$ws.Range("V5").Value = 19
$ws.Range("W5").Value = 120
$ws.Range("X5").Value = 55
$ws.Range("Y5").Value = 45
$ws.Range("Z5").Value = 17
$ws.Range("AA5").Value = 9
$ws.Range("AB5").Value = 15.5
$ws.Range("AC5").Value = 55
$ws.Range("AD5").Value = 300
$ws.Range("AE5").Value = 7.9
$ws.Range("AF5").Value = 6.5
$ws.Range("AG5").Value = 7.4
$ws.Range("AH5").Value = 7.5
$ws.Range("AI5").Value = 8.5
$ws.Range("AJ5").Value = 18

# Row 6 - odds newly populated for Tammeka vs Tallinna Kalev
$ws.Range("G6").Value = 2.27
$ws.Range("H6").Value = 3.45
$ws.Range("I6").Value = 2.65
$ws.Range("T6").Value = 10.25
$ws.Range("U6").Value = 12.5
$ws.Range("V6").Value = 7.9
$ws.Range("W6").Value = 21
$ws.Range("X6").Value = 13.5
$ws.Range("Y6").Value = 15.5
$ws.Range("Z6").Value = 15
$ws.Range("AA6").Value = 6.4
$ws.Range("AB6").Value = 9.25
$ws.Range("AC6").Value = 26
$ws.Range("AD6").Value = 120
$ws.Range("AE6").Value = 10.75
$ws.Range("AF6").Value = 14
$ws.Range("AG6").Value = 8.5
$ws.Range("AH6").Value = 26
$ws.Range("AI6").Value = 16
$ws.Range("AJ6").Value = 17.5

# Row 7 - odds newly populated for Narva vs Harju JK Laagri
$ws.Range("G7").Value = 1.62
$ws.Range("H7").Value = 3.8
$ws.Range("I7").Value = 4.3
$ws.Range("T7").Value = 9.25
$ws.Range("U7").Value = 9.25
$ws.Range("V7").Value = 7.1
$ws.Range("W7").Value = 12.5
$ws.Range("X7").Value = 9.5
$ws.Range("Y7").Value = 13.5
$ws.Range("Z7").Value = 16.5
$ws.Range("AA7").Value = 7.1
$ws.Range("AB7").Value = 10.25
$ws.Range("AC7").Value = 29
$ws.Range("AD7").Value = 150
$ws.Range("AE7").Value = 14.5
$ws.Range("AF7").Value = 24
$ws.Range("AG7").Value = 12
$ws.Range("AH7").Value = 55
$ws.Range("AI7").Value = 28
$ws.Range("AJ7").Value = 25

# Row 8 - odds updated for Braunschweig vs Saarbrucken
$ws.Range("G8").Value = 1.98
$ws.Range("H8").Value = 3.55
$ws.Range("L8").Value = 1.3
$ws.Range("M8").Value = 3.6
$ws.Range("N8").Value = 1.88
$ws.Range("O8").Value = 1.93
$ws.Range("Q8").Value = 2.82
$ws.Range("R8").Value = 1.72
$ws.Range("S8").Value = 2.05
$ws.Range("T8").Value = 7.3
$ws.Range("U8").Value = 10
$ws.Range("W8").Value = 18.5
$ws.Range("X8").Value = 16.5
$ws.Range("Y8").Value = 28
$ws.Range("AA8").Value = 7
$ws.Range("AB8").Value = 14.5
$ws.Range("AC8").Value = 65
$ws.Range("AD8").Value = 500
$ws.Range("AE8").Value = 11.25
$ws.Range("AF8").Value = 23
$ws.Range("AG8").Value = 13
$ws.Range("AH8").Value = 60
$ws.Range("AI8").Value = 35
$ws.Range("AJ8").Value = 40

# Row 9 - odds updated for Daegu vs Jeonbuk
$ws.Range("T9").Value = 15

# Row 10 - odds updated for Daejeon vs Pohang
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 3.25
$ws.Range("K10").Value = 10
$ws.Range("L10").Value = 1.33
$ws.Range("M10").Value = 3.25
$ws.Range("N10").Value = 2.05
$ws.Range("O10").Value = 1.75
$ws.Range("U10").Value = 11
$ws.Range("V10").Value = 9.5
$ws.Range("AA10").Value = 6
$ws.Range("AF10").Value = 15

# Row 11 - odds updated for Suwon FC vs Jeju SK
$ws.Range("K11").Value = 10
